$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TLHP")

# Insert a new row at row 3 (shifts existing rows 3..50 down to 4..51)
$ws.Rows(3).Insert()

# Populate the new row's values
$ws.Range("B3").Value = "User Authentication "
$ws.Range("C3").Value = "done"

# Hide the new row (it matches status "done" which is filtered out)
$ws.Rows(3).Hidden = $true

# Restore the active selection as left by the editor
$ws.Range("B15").Select()
